# "Generate Report for Handoff"
# The b61a95ae-9908-48a1-8db0-f601d2fb8268.md file has been handed off for
# localization: its status flips from "In Translation" to "Ready for
# handoff" on the Overview sheet (both locale columns) as well as on each
# per-locale detail sheet, the handoff priority becomes "mt", and the
# handoff timestamps are refreshed.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 == b61a95ae-9908-48a1-8db0-f601d2fb8268.md ---
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-04 12:17:47"

# --- zh-cn sheet: row 3 == b61a95ae-9908-48a1-8db0-f601d2fb8268.md ---
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-09-04 12:17:43"

# --- de-de sheet: row 3 == b61a95ae-9908-48a1-8db0-f601d2fb8268.md ---
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-09-04 12:17:47"

# --- Column widths widen slightly to fit the longer "Ready for handoff" /
#     "Latest Handoff Datetime" text (target stored width ~17.216; COM
#     quantizes ColumnWidth to a 1/6-pixel grid so 16.3 is the closest
#     achievable input, landing on 17.1667). ---
$wsOverview.Range("E1").ColumnWidth = 16.3
$wsOverview.Range("F1").ColumnWidth = 16.3
$wsZhCn.Range("C1").ColumnWidth = 16.3
$wsDeDe.Range("C1").ColumnWidth = 16.3
